# ToDo.xlsx update:
#  - Two finished tasks in B6/B7 are replaced by two new tasks
#    ("test-that implementieren" / "renv implementieren")
#  - The two-row block "Zusammenfassung der Transformation Monitors" /
#    "Messaging Objekt" (rows 12-13, plus the blank spacer row 14) is
#    removed entirely, shifting all rows below it up.
#  - Selection moves to B18 (the blank row right after the last task).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace completed task descriptions with the new tasks.
$ws.Range("B6").Value = "test-that implementieren"
$ws.Range("B7").Value = "renv implementieren"

# Remove the now-irrelevant "Zusammenfassung der Transformation Monitors" /
# "Messaging Objekt" rows (and the blank row following them), shifting
# everything below up by three rows.
$ws.Range("B12:B14").EntireRow.Delete() | Out-Null

# Move the active selection to match the saved view state.
$ws.Range("B18").Select() | Out-Null
